$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 29 (RS / Serbia) - Excel shifts all subsequent rows up by one
# and automatically re-points the SUMIF formula range (G2:G33/C2:C33 -> G2:G32/C2:C32).
$ws.Rows.Item(29).Delete()

# Re-apply the sort so the sortState/sortCondition range shrinks to match
# the new data extent (A2:A32 instead of A2:A33).
$ws.Sort.SortFields.Add2($ws.Range("A2:A32"))
$ws.Sort.SetRange($ws.Range("A2:A32"))
$ws.Sort.Apply()

# Set the active selection cell, matching the final saved state.
$ws.Range("F24").Select()
